# The "parameters" sheet (sheetId 1) is the active/tab-selected sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L (rows 41-74) builds up PHP assignment lines such as:
#   $member->id = $data->id;
# The author renamed the local variable from $member to $item, e.g.:
#   $item->id = $data->id;
# Column C on the corresponding header row (C1, C2, ... C34) holds the field name
# referenced by each of these formulas, mirroring the existing pattern used by the
# sibling columns B/D/H (e.g. "$this->"&C1&" = $row['"&C1&"'];").
for ($i = 0; $i -lt 34; $i++) {
    $row = 41 + $i
    $cRef = "C" + ($i + 1)
    $ws.Range("L$row").Formula = '="$item->"&' + $cRef + '&" = $data->"&' + $cRef + '&";"'
}

# Leave the selection on L1, matching where the edit was finished.
$ws.Range("L1").Select()
